# Update the handback-status report with refreshed timestamps / status
# as produced by a later run of the "Generate Report for Handback" job.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for 594caa70-f1de-4bc8-9f88-7e20074aafd6
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-29 16:17:58"
$wsOverview.Range("G5").Value = "2016-08-29 16:17:58"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column for the 594caa70 row changed from "ht" to "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H3").Value = "2016-08-29 16:17:52"
$wsZhCn.Range("H5").Value = "2016-08-29 16:17:52"
# Correspond Handback DateTime
$wsZhCn.Range("K3").Value = "2016-08-29 16:18:18"
$wsZhCn.Range("K5").Value = "2016-08-29 16:18:18"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Correspond Handoff Datetime
$wsDeDe.Range("H3").Value = "2016-08-29 16:17:58"
$wsDeDe.Range("H5").Value = "2016-08-29 16:17:58"
# Priority column for the 594caa70 row changed from "ht" to "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
# Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-08-29 16:18:26"
$wsDeDe.Range("K5").Value = "2016-08-29 16:18:26"
